# Updated cryptos list on Tue Oct 17 13:00:04 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "28.372.78"
Set-TextCell "E2" "  +2.90%  "
Set-TextCell "D3" "1.581.82"
Set-TextCell "E3" "  +0.50%  "
Set-TextCell "D5" "212.17"
Set-TextCell "E5" "  +0.43%  "
Set-TextCell "E6" "  -0.08%  "
Set-TextCell "E7" "  +1.20%  "
Set-TextCell "B8" "OKB"
Set-TextCell "C8" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D8" "46.59"
Set-TextCell "E8" "  +7.10%  "
Set-TextCell "B9" "Solana"
Set-TextCell "C9" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell "D9" "23.99"
Set-TextCell "E9" "  +3.07%  "
Set-TextCell "B10" "Cardano"
Set-TextCell "C10" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextCell "D10" "0.250"
Set-TextCell "E10" "  -0.50%  "
Set-TextCell "B11" "Dogecoin"
Set-TextCell "C11" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextCell "D11" "0.0597"
Set-TextCell "E11" "  -0.16%  "
Set-TextCell "B12" "TRON"
Set-TextCell "C12" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D12" "0.0883"
Set-TextCell "E12" "  +1.11%  "
Set-TextCell "B13" "WrappedliquidstakedEther2.0"
Set-TextCell "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell "D13" "1.809.53"
Set-TextCell "E13" "  +0.62%  "
Set-TextCell "B14" "WrappedEther"
Set-TextCell "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D14" "1.583.32"
Set-TextCell "E14" "  +1.13%  "
Set-TextCell "B15" "Polygon"
Set-TextCell "C15" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D15" "0.527"
Set-TextCell "E15" "  +0.91%  "
Set-TextCell "B16" "Polkadot"
Set-TextCell "C16" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D16" "3.72"
Set-TextCell "E16" "  -1.33%  "
Set-TextCell "B17" "WrappedBTC"
Set-TextCell "C17" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell "D17" "28.444.86"
Set-TextCell "E17" "  +3.38%  "
Set-TextCell "B18" "Litecoin"
Set-TextCell "C18" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextCell "D18" "62.62"
Set-TextCell "E18" "  -1.10%  "
Set-TextCell "B19" "BitcoinCash"
Set-TextCell "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D19" "229.57"
Set-TextCell "E19" "  -0.36%  "
Set-TextCell "E20" "  -0.78%  "
Set-TextCell "B21" "Chainlink"
Set-TextCell "C21" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D21" "7.42"
Set-TextCell "E21" "  -1.46%  "
Set-TextCell "B22" "Dai"
Set-TextCell "C22" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D22" "1.00"
Set-TextCell "E22" "  +0.95%  "
Set-TextCell "B23" "Uniswap"
Set-TextCell "C23" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D23" "3.95"
Set-TextCell "E23" "  -4.09%  "
Set-TextCell "B24" "Avalanche"
Set-TextCell "C24" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D24" "9.24"
Set-TextCell "E24" "  -2.29%  "
Set-TextCell "B25" "Toncoin"
Set-TextCell "C25" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D25" "2.02"
Set-TextCell "E25" "  +3.08%  "
Set-TextCell "B26" "Monero"
Set-TextCell "C26" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell "D26" "151.36"
Set-TextCell "E26" "  +0.72%  "
Set-TextCell "B27" "EthereumClassic"
Set-TextCell "C27" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D27" "15.11"
Set-TextCell "E27" "  -1.61%  "
Set-TextCell "B28" "Cosmos"
Set-TextCell "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D28" "6.51"
Set-TextCell "E28" "  -1.46%  "
Set-TextCell "B29" "Stellar"
Set-TextCell "C29" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D29" "0.105"
Set-TextCell "E29" "  -1.79%  "
Set-TextCell "B30" "BinanceUSD"
Set-TextCell "C30" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D30" "1.00"
Set-TextCell "E30" "  +1.09%  "
Set-TextCell "B31" "PancakeSwap"
Set-TextCell "C31" "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell "D31" "1.12"
Set-TextCell "E31" "  -1.33%  "
Set-TextCell "B32" "Hedera"
Set-TextCell "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D32" "0.0467"
Set-TextCell "E32" "  -1.50%  "
Set-TextCell "B33" "Filecoin"
Set-TextCell "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D33" "3.23"
Set-TextCell "E33" "  -0.66%  "
Set-TextCell "B34" "InternetComputer(DFINITY)"
Set-TextCell "C34" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D34" "3.14"
Set-TextCell "E34" "  +0.30%  "
Set-TextCell "B35" "Maker"
Set-TextCell "C35" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D35" "1.388.49"
Set-TextCell "E35" "  -5.01%  "
Set-TextCell "B36" "LidoDAOToken"
Set-TextCell "C36" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D36" "1.57"
Set-TextCell "E36" "  -2.49%  "
Set-TextCell "B37" "TrustWalletToken"
Set-TextCell "C37" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D37" "1.01"
Set-TextCell "E37" "  -5.37%  "
Set-TextCell "B38" "HuobiToken"
Set-TextCell "C38" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D38" "2.35"
Set-TextCell "E38" "  +1.03%  "
Set-TextCell "B39" "MXToken"
Set-TextCell "C39" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D39" "2.61"
Set-TextCell "E39" "  +9.61%  "
Set-TextCell "B40" "VeChain"
Set-TextCell "C40" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D40" "0.0166"
Set-TextCell "E40" "  -1.32%  "
Set-TextCell "B41" "ImmutableX"
Set-TextCell "C41" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell "D41" "0.535"
Set-TextCell "E41" "  -1.15%  "
Set-TextCell "B42" "ARBITRUM"
Set-TextCell "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D42" "0.805"
Set-TextCell "E42" "  -1.29%  "
Set-TextCell "B43" "PaxDollar"
Set-TextCell "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D43" "1.00"
Set-TextCell "E43" "  +1.05%  "
Set-TextCell "E44" "  +1.26%  "
Set-TextCell "B45" "FraxShare"
Set-TextCell "C45" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D45" "5.59"
Set-TextCell "E45" "  -1.64%  "
Set-TextCell "B46" "WEMIXToken"
Set-TextCell "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D46" "0.982"
Set-TextCell "E46" "  +1.00%  "
Set-TextCell "B47" "Aave"
Set-TextCell "C47" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D47" "62.55"
Set-TextCell "E47" "  -2.69%  "
Set-TextCell "B48" "RocketPoolETH"
Set-TextCell "C48" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell "D48" "1.718.97"
Set-TextCell "E48" "  +0.65%  "
Set-TextCell "D49" "86.30"
Set-TextCell "E49" "  -0.63%  "
Set-TextCell "B50" "Cronos"
Set-TextCell "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D50" "0.0519"
Set-TextCell "E50" "  -0.93%  "
Set-TextCell "B51" "BabyDogeCoin"
Set-TextCell "C51" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D51" "0.0₇0999"
Set-TextCell "E51" "  -0.01%  "
